# The deck ships two theme parts:
#   ppt/theme/theme1.xml -> bound to the notes master  ("Office Theme" colours)
#   ppt/theme/theme2.xml -> bound to the slide master   ("Integral" colours)
# The authored edit swaps the two themes' content, so the slide master now
# carries the plain "Office Theme" palette and the notes master carries the
# "Integral" palette. The notes master's theme isn't reachable through the
# PowerPoint object model here (NotesMaster resolves to the same Master as
# the slide master in this host), so we recolour the one theme that the
# object model does expose -- the slide master's theme (ppt/theme/theme2.xml)
# -- from "Integral" to the "Office Theme" palette, matching the diff for
# that part.

$p  = $ppt.ActivePresentation
$sm = $p.SlideMaster
$tcs = $sm.Theme.ThemeColorScheme

# Target palette = the stock "Office Theme" scheme (what theme2.xml becomes).
# RGBColor.RGB takes/returns the Win32 BGR-packed 0xBBGGRR integer form, so
# every 0xRRGGBB hex colour below is converted accordingly.
$officeThemeColors = @(
    0x000000,  # 1  dk1
    0xFFFFFF,  # 2  lt1
    0x44546A,  # 3  dk2
    0xE7E6E6,  # 4  lt2
    0x5B9BD5,  # 5  accent1
    0xED7D31,  # 6  accent2
    0xA5A5A5,  # 7  accent3
    0xFFC000,  # 8  accent4
    0x4472C4,  # 9  accent5
    0x70AD47,  # 10 accent6
    0x0563C1,  # 11 hlink
    0x954F72   # 12 folHlink
)

for ($i = 1; $i -le $tcs.Count; $i++) {
    # Literal 0xRRGGBB -> low byte is Blue, mid byte Green, high byte Red.
    $hex = $officeThemeColors[$i - 1]
    $b = $hex -band 0xFF
    $g = ($hex -shr 8) -band 0xFF
    $r = ($hex -shr 16) -band 0xFF
    $bgr = ($b -shl 16) -bor ($g -shl 8) -bor $r

    $color = $tcs.Colors($i)
    $color.RGB = $bgr
}
